# "Use improved tooth profile" - add a "Pulley Design" worksheet with belt /
# pulley tooth-profile reference data (Gates & Optibelt datasheets, pitch /
# outside diameter interpolation) after the existing Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# New worksheet, positioned after Sheet1, becomes the active/selected sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Pulley Design"

# ---------------------------------------------------------------------------
# Header row (row 1) - bold
# ---------------------------------------------------------------------------
$ws2.Range("C1").Value = "Name"
$ws2.Range("D1").Value = "URL"
$ws2.Range("E1").Value = "Pages"
$ws2.Range("M1").Value = "Source"
$ws2.Range("C1:E1").Font.Bold = $true
$ws2.Range("M1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Reference list (rows 2-4): [n] marker, document name + hyperlinked URL,
# and relevant page(s)
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = "[1]"
$ws2.Range("C2").Value = "Gates Powergrip Drive Design Manual"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://www.gates.com/content/dam/gates/home/knowledge-center/resource-library/catalogs/powergripdrivedesignmanual_17195_2014.pdf") | Out-Null
$ws2.Range("E2").Value = "131, 138"

$ws2.Range("B3").Value = "[2]"
$ws2.Range("C3").Value = "Technical Data Sheet - optibelt OMEGA HP 8M"
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://www.optibelt.com/fileadmin/pdf/datenblaetter/Technical-Data-Sheet-optibelt-OMEGA-HP-8M.pdf") | Out-Null
$ws2.Range("E3").Value = 1

$ws2.Range("B4").Value = "[3]"
$ws2.Range("C4").Value = "Optibelt-TM-Rubber-Timing-Belt-Drives"
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://www.optibelt.com/fileadmin/pdf/produkte/zahnriemen-gummi/Optibelt-TM-Rubber-Timing-Belt-Drives.pdf") | Out-Null
$ws2.Range("E4").Value = 75

# ---------------------------------------------------------------------------
# Teeth / pitch-diameter / outside-diameter lookup table (J2:M6)
# ---------------------------------------------------------------------------
$ws2.Range("J2").Value = "Teeth"
$ws2.Range("K2").Value = 39
$ws2.Range("L2").Value = 63

$ws2.Range("J3").Value = "Pitch Dia"
$ws2.Range("K3").Value = 99.31
$ws2.Range("L3").Value = 160.43
$ws2.Range("M3").Value = "[3]"

$ws2.Range("J4").Value = "Outside Dia"
$ws2.Range("K4").Value = 97.94
$ws2.Range("L4").Value = 159.06
$ws2.Range("M4").Value = "[3]"

$ws2.Range("J5").Value = "Pitch Dia (Alt)"
$ws2.Range("K5").Formula = '=3.91*25.4'
$ws2.Range("L5").Formula = '=S4'
$ws2.Range("M5").Value = "[1]"

$ws2.Range("J6").Value = "Outside Dia (Alt)"
$ws2.Range("K6").Formula = '=3.856*25.4'
$ws2.Range("L6").Formula = '=T4'
$ws2.Range("M6").Value = "[1]"

$ws2.Range("J2:J6").Font.Bold = $true

# ---------------------------------------------------------------------------
# Interpolation helper table (O3:T4) - 56/63/64-tooth pulley pitch & outside
# diameters (inches), with mm conversions
# ---------------------------------------------------------------------------
$ws2.Range("O3").Value = 56
$ws2.Range("P3").Value = 5.614
$ws2.Range("Q3").Value = 5.56

$ws2.Range("O4").Value = 63
$ws2.Range("P4").Formula = '=(P5-P3)/($O5-$O3)*($O4-$O3)+P3'
$ws2.Range("Q4").Formula = '=(Q5-Q3)/($O5-$O3)*($O4-$O3)+Q3'
$ws2.Range("S4").Formula = '=P4*25.4'
$ws2.Range("T4").Formula = '=Q4*25.4'

$ws2.Range("O5").Value = 64
$ws2.Range("P5").Value = 6.416
$ws2.Range("Q5").Value = 6.362

$ws2.Range("O6").Font.Bold = $true
$ws2.Range("O3:O5").Font.Bold = $true

# mm-precision number formats used by the S/T helper columns (reuse the
# existing "0.00" style already present in the workbook)
$ws2.Range("S4:T4").NumberFormat = "0.00"
$ws2.Range("K5:L6").NumberFormat = "0.00"

# decimal-precision number formats for the O:T interpolation block - applied
# in this order so the generated numFmtIds line up (164=0.00000,
# 165=0.000, 166=0.000000, 167=0.00000000, 168=0.0000000000)
$ws2.Range("P11").NumberFormat = "0.00000"
$ws2.Range("P3:Q5").NumberFormat = "0.000"
$ws2.Range("Q11:R11").NumberFormat = "0.000000"
$ws2.Range("P15:P17").NumberFormat = "0.00000000"
$ws2.Range("R16:R17").NumberFormat = "0.00000000"
$ws2.Range("Q15:Q17").NumberFormat = "0.0000000000"

$ws2.Range("O15:O17").Font.Bold = $true

# ---------------------------------------------------------------------------
# Column widths (approximate - engine stores widths at coarser granularity
# than native Excel)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 32.89
$ws2.Columns.Item(10).ColumnWidth = 12.44
$ws2.Columns.Item(16).ColumnWidth = 12.55
$ws2.Range("Q16:R17").ColumnWidth = 13.66

# ---------------------------------------------------------------------------
# Selection / active sheet
# ---------------------------------------------------------------------------
$ws2.Range("M15").Select()
$ws2.Select()

$wb.Application.Calculate()
